# Apply scheduled market-data refresh to price/profit columns (H-N)
# across the ALC, ARM, BSM, CRP, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1556.6666
$ws.Range("I43").Value = 1722.5
$ws.Range("J43").Value = 1225
$ws.Range("K43").Value = 1722.5
$ws.Range("L43").Value = 1225
$ws.Range("M43").Value = -1653.5
$ws.Range("N43").Value = -1363

$ws.Range("H80").Value = 1422.5
$ws.Range("I80").Value = 146.66667
$ws.Range("J80").Value = 5250
$ws.Range("K80").Value = 440.00001
$ws.Range("L80").Value = 15750
$ws.Range("M80").Value = 557.99999
$ws.Range("N80").Value = -17746

$ws.Range("H83").Value = 1422.5
$ws.Range("I83").Value = 146.66667
$ws.Range("J83").Value = 5250
$ws.Range("K83").Value = 1320.00003
$ws.Range("L83").Value = 47250
$ws.Range("M83").Value = 3671.99997
$ws.Range("N83").Value = -57234

$ws.Range("H111").Value = 83335210
$ws.Range("I111").Value = 191.8
$ws.Range("J111").Value = 142860210
$ws.Range("K111").Value = 575.4000000000001
$ws.Range("L111").Value = 428580630
$ws.Range("M111").Value = 2491.6
$ws.Range("N111").Value = -428586764

$ws.Range("H137").Value = 1273.2
$ws.Range("I137").Value = 1027.9459
$ws.Range("J137").Value = 4298
$ws.Range("K137").Value = 3083.8377
$ws.Range("L137").Value = 12894
$ws.Range("M137").Value = -533.8377
$ws.Range("N137").Value = -17994

$ws.Range("H138").Value = 2364.4614
$ws.Range("I138").Value = 1183.9269
$ws.Range("J138").Value = 3672.6216
$ws.Range("K138").Value = 3551.7807
$ws.Range("L138").Value = 11017.8648
$ws.Range("M138").Value = 1588.2193
$ws.Range("N138").Value = -21297.8648

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1836.375
$ws.Range("I61").Value = 1891.5714
$ws.Range("J61").Value = 1450
$ws.Range("K61").Value = 1891.5714
$ws.Range("L61").Value = 1450
$ws.Range("M61").Value = -1679.5714
$ws.Range("N61").Value = -1874

$ws.Range("H74").Value = 2359.5
$ws.Range("I74").Value = 2900.4
$ws.Range("J74").Value = 1458
$ws.Range("K74").Value = 2900.4
$ws.Range("L74").Value = 1458
$ws.Range("M74").Value = -2026.4
$ws.Range("N74").Value = -3206

$ws.Range("H77").Value = 2359.5
$ws.Range("I77").Value = 2900.4
$ws.Range("J77").Value = 1458
$ws.Range("K77").Value = 14502
$ws.Range("L77").Value = 7290
$ws.Range("M77").Value = -10134
$ws.Range("N77").Value = -16026

$ws.Range("H132").Value = 2915.2632
$ws.Range("I132").Value = 1426.7059
$ws.Range("J132").Value = 4120.2856
$ws.Range("K132").Value = 4280.1177
$ws.Range("L132").Value = 12360.8568
$ws.Range("M132").Value = -1750.1177
$ws.Range("N132").Value = -17420.8568

$ws.Range("H136").Value = 1836.375
$ws.Range("I136").Value = 1891.5714
$ws.Range("J136").Value = 1450
$ws.Range("K136").Value = 5674.7142
$ws.Range("L136").Value = 4350
$ws.Range("M136").Value = -3124.7142
$ws.Range("N136").Value = -9450

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H120").Value = 20000
$ws.Range("J120").Value = 20000
$ws.Range("L120").Value = 20000
$ws.Range("N120").Value = -29676

$ws.Range("H134").Value = 2372.5173
$ws.Range("I134").Value = 1391.8235
$ws.Range("J134").Value = 3761.8333
$ws.Range("K134").Value = 4175.470499999999
$ws.Range("L134").Value = 11285.4999
$ws.Range("M134").Value = -1640.470499999999
$ws.Range("N134").Value = -16355.4999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2778.2307
$ws.Range("I31").Value = 2411.652
$ws.Range("J31").Value = 3305.1875
$ws.Range("K31").Value = 2411.652
$ws.Range("L31").Value = 3305.1875
$ws.Range("M31").Value = -2116.652
$ws.Range("N31").Value = -3895.1875

$ws.Range("H34").Value = 2778.2307
$ws.Range("I34").Value = 2411.652
$ws.Range("J34").Value = 3305.1875
$ws.Range("K34").Value = 2411.652
$ws.Range("L34").Value = 3305.1875
$ws.Range("M34").Value = -2209.652
$ws.Range("N34").Value = -3709.1875

$ws.Range("H58").Value = 2587.7585
$ws.Range("I58").Value = 1790.6364
$ws.Range("J58").Value = 3074.889
$ws.Range("K58").Value = 1790.6364
$ws.Range("L58").Value = 3074.889
$ws.Range("M58").Value = -1587.6364
$ws.Range("N58").Value = -3480.889

$ws.Range("H132").Value = 5226.25
$ws.Range("I132").Value = 4453
$ws.Range("J132").Value = 5999.5
$ws.Range("K132").Value = 13359
$ws.Range("L132").Value = 17998.5
$ws.Range("M132").Value = -10829
$ws.Range("N132").Value = -23058.5

$ws.Range("H134").Value = 2924.25
$ws.Range("I134").Value = 2686.5293
$ws.Range("J134").Value = 4271.3335
$ws.Range("K134").Value = 8059.5879
$ws.Range("L134").Value = 12814.0005
$ws.Range("M134").Value = -5524.5879
$ws.Range("N134").Value = -17884.0005

$ws.Range("H136").Value = 2587.7585
$ws.Range("I136").Value = 1790.6364
$ws.Range("J136").Value = 3074.889
$ws.Range("K136").Value = 5371.9092
$ws.Range("L136").Value = 9224.667000000001
$ws.Range("M136").Value = -2821.9092
$ws.Range("N136").Value = -14324.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2773.2
$ws.Range("I132").Value = 2049.6191
$ws.Range("J132").Value = 4461.5557
$ws.Range("K132").Value = 6148.8573
$ws.Range("L132").Value = 13384.6671
$ws.Range("M132").Value = -3618.8573
$ws.Range("N132").Value = -18444.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6473.1924
$ws.Range("I132").Value = 7802.6177
$ws.Range("J132").Value = 3962.0557
$ws.Range("K132").Value = 23407.8531
$ws.Range("L132").Value = 11886.1671
$ws.Range("M132").Value = -20877.8531
$ws.Range("N132").Value = -16946.1671

$ws.Range("H136").Value = 3262.8286
$ws.Range("I136").Value = 3025.75
$ws.Range("J136").Value = 3578.9333
$ws.Range("K136").Value = 9077.25
$ws.Range("L136").Value = 10736.7999
$ws.Range("M136").Value = -6527.25
$ws.Range("N136").Value = -15836.7999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2517.5715
$ws.Range("I132").Value = 2110.0588
$ws.Range("J132").Value = 4249.5
$ws.Range("K132").Value = 6330.176399999999
$ws.Range("L132").Value = 12748.5
$ws.Range("M132").Value = -3800.176399999999
$ws.Range("N132").Value = -17808.5

$ws.Range("H136").Value = 13892697
$ws.Range("I136").Value = 22223222
$ws.Range("J136").Value = 8490
$ws.Range("K136").Value = 66669666
$ws.Range("L136").Value = 25470
$ws.Range("M136").Value = -66667116
$ws.Range("N136").Value = -30570
